# Applies the two-digit / one-digit division answer-key updates described in the diff.
# Each change is targeted at a specific (row, column) cell of the single table in the
# document, so the one duplicated value ("62÷4=15, 2", which appears twice in the original)
# cannot be mismatched the way a blind Find/Replace-All could mismatch it.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1
$cell = $t.Rows(1).Cells(1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "30÷5=6, 0") {
    throw "Unexpected text in row 1 cell 1: $($cell.Range.Text)"
}
$cell.Range.Text = "39÷6=6, 3"

$cell = $t.Rows(1).Cells(2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "96÷9=10, 6") {
    throw "Unexpected text in row 1 cell 2: $($cell.Range.Text)"
}
$cell.Range.Text = "38÷9=4, 2"

$cell = $t.Rows(1).Cells(3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "75÷7=10, 5") {
    throw "Unexpected text in row 1 cell 3: $($cell.Range.Text)"
}
$cell.Range.Text = "39÷8=4, 7"

$cell = $t.Rows(1).Cells(4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "28÷6=4, 4") {
    throw "Unexpected text in row 1 cell 4: $($cell.Range.Text)"
}
$cell.Range.Text = "59÷8=7, 3"

$cell = $t.Rows(1).Cells(5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "82÷2=41, 0") {
    throw "Unexpected text in row 1 cell 5: $($cell.Range.Text)"
}
$cell.Range.Text = "58÷4=14, 2"

# Row 5
$cell = $t.Rows(5).Cells(1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "63÷9=7, 0") {
    throw "Unexpected text in row 5 cell 1: $($cell.Range.Text)"
}
$cell.Range.Text = "45÷9=5, 0"

$cell = $t.Rows(5).Cells(2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "55÷5=11, 0") {
    throw "Unexpected text in row 5 cell 2: $($cell.Range.Text)"
}
$cell.Range.Text = "30÷3=10, 0"

$cell = $t.Rows(5).Cells(3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "39÷8=4, 7") {
    throw "Unexpected text in row 5 cell 3: $($cell.Range.Text)"
}
$cell.Range.Text = "30÷9=3, 3"

$cell = $t.Rows(5).Cells(4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "64÷4=16, 0") {
    throw "Unexpected text in row 5 cell 4: $($cell.Range.Text)"
}
$cell.Range.Text = "25÷2=12, 1"

$cell = $t.Rows(5).Cells(5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "11÷5=2, 1") {
    throw "Unexpected text in row 5 cell 5: $($cell.Range.Text)"
}
$cell.Range.Text = "93÷7=13, 2"

# Row 9
$cell = $t.Rows(9).Cells(1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "71÷7=10, 1") {
    throw "Unexpected text in row 9 cell 1: $($cell.Range.Text)"
}
$cell.Range.Text = "41÷7=5, 6"

$cell = $t.Rows(9).Cells(2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "47÷3=15, 2") {
    throw "Unexpected text in row 9 cell 2: $($cell.Range.Text)"
}
$cell.Range.Text = "93÷4=23, 1"

$cell = $t.Rows(9).Cells(3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "62÷4=15, 2") {
    throw "Unexpected text in row 9 cell 3: $($cell.Range.Text)"
}
$cell.Range.Text = "53÷5=10, 3"

$cell = $t.Rows(9).Cells(4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "20÷2=10, 0") {
    throw "Unexpected text in row 9 cell 4: $($cell.Range.Text)"
}
$cell.Range.Text = "37÷9=4, 1"

$cell = $t.Rows(9).Cells(5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "29÷5=5, 4") {
    throw "Unexpected text in row 9 cell 5: $($cell.Range.Text)"
}
$cell.Range.Text = "52÷3=17, 1"

# Row 13
$cell = $t.Rows(13).Cells(1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "62÷4=15, 2") {
    throw "Unexpected text in row 13 cell 1: $($cell.Range.Text)"
}
$cell.Range.Text = "26÷6=4, 2"

$cell = $t.Rows(13).Cells(2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "31÷7=4, 3") {
    throw "Unexpected text in row 13 cell 2: $($cell.Range.Text)"
}
$cell.Range.Text = "35÷3=11, 2"

$cell = $t.Rows(13).Cells(3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "77÷7=11, 0") {
    throw "Unexpected text in row 13 cell 3: $($cell.Range.Text)"
}
$cell.Range.Text = "83÷7=11, 6"

$cell = $t.Rows(13).Cells(4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "31÷5=6, 1") {
    throw "Unexpected text in row 13 cell 4: $($cell.Range.Text)"
}
$cell.Range.Text = "73÷5=14, 3"

$cell = $t.Rows(13).Cells(5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "11÷6=1, 5") {
    throw "Unexpected text in row 13 cell 5: $($cell.Range.Text)"
}
$cell.Range.Text = "20÷9=2, 2"

# Row 17
$cell = $t.Rows(17).Cells(1)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "53÷7=7, 4") {
    throw "Unexpected text in row 17 cell 1: $($cell.Range.Text)"
}
$cell.Range.Text = "25÷5=5, 0"

$cell = $t.Rows(17).Cells(2)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "84÷3=28, 0") {
    throw "Unexpected text in row 17 cell 2: $($cell.Range.Text)"
}
$cell.Range.Text = "77÷8=9, 5"

$cell = $t.Rows(17).Cells(3)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "49÷8=6, 1") {
    throw "Unexpected text in row 17 cell 3: $($cell.Range.Text)"
}
$cell.Range.Text = "90÷6=15, 0"

$cell = $t.Rows(17).Cells(4)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "90÷6=15, 0") {
    throw "Unexpected text in row 17 cell 4: $($cell.Range.Text)"
}
$cell.Range.Text = "83÷3=27, 2"

$cell = $t.Rows(17).Cells(5)
if ($cell.Range.Text.TrimEnd([char]13,[char]7) -ne "82÷5=16, 2") {
    throw "Unexpected text in row 17 cell 5: $($cell.Range.Text)"
}
$cell.Range.Text = "69÷6=11, 3"

